$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "HoaiBao"
$ws.Range("B3").Value = "123456aA@"
